# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at the top of the data
# (row 2), pushing all existing rows down by one. The new row reuses the
# same market/category/variety/quality/unit/origin/kg metadata as the
# (old) first data row, but carries its own date and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 2 (shifts rows 2-26 -> 3-27).
$ws.Rows.Item(2).Insert()

# The default Insert() copies the header row's (bold/bordered) formatting
# onto the new row; strip that back to the plain "data row" look used by
# every other record, then restore the date number format on column D.
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with this week's reading.
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44643
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 100112022
$ws.Cells.Item(2, 7).Value = "Arveja Verde"
$ws.Cells.Item(2, 8).Value = "Perfection"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 90
$ws.Cells.Item(2, 11).Value = 25000
$ws.Cells.Item(2, 12).Value = 26000
$ws.Cells.Item(2, 13).Value = 25444
$ws.Cells.Item(2, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(2, 15).Value = "Carahue"
$ws.Cells.Item(2, 16).Value = 1018
$ws.Cells.Item(2, 17).Value = 25
$ws.Cells.Item(2, 18).Value = "Hortaliza"
